$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$intent = 'Não sei'
$sentences = @(
    'Quais os meus compromissos de hoje?',
    'Quanto foi o jogo do Palmeiras?',
    'Foxbot, como faço para chegar em casa?',
    'Hoje tem jogo do Barcelona?',
    'Mostre uma foto do Neymar',
    'mostre uma receita de caipirinha',
    'Quem ganhou o último BBB?',
    'Mostre a meus compromissos para a semana',
    'Consulte o número de mortes por coronavírus hoje',
    'O Palmeiras tem mundial?',
    'Quando foi o último título do São Paulo?',
    'Conte uma piada',
    'Quem é a Gabi da AE4?',
    'O Matheus é um autômato?',
    'Como faço para baixar RAM? ',
    'Marque uma consulta com o meu médico para semana que vem.',
    'Você me ama?',
    'Foxbot, qual a melhor rota para chegar em casa?',
    'Como faço para aumentar meu pênis?',
    'Por que o Ribamar não é da seleção?',
    'Receita de miojo',
    'De que cookie você gosta?',
    'Você acredita em Deus?',
    'Ainda tenho compromissos hoje?',
    'Foxbot, o Acre existe?',
    'Foxbot, toque Grupo Menos é Mais',
    'Mostre fotos do carnaval passado',
    'Vacina transforma em jacaré?',
    'O que é um mamaco?',
    'Marque uma reunião com o estagiário para amanhã as 17 horas',
    'Por que não escovar os dentes com Hipoglós?',
    'Foxbot, bora fumar um?',
    'Foxbot, toque Thiaguinho',
    'Foxbot, coloque um alarme para amanhã às 7h30',
    'Me acorde amanhã às 9 horas',
    'Malboro é um bom pré treino?',
    'Foxbot, ligue o cronômetro',
    'Foxbot, ligue para a maravilhosa',
    'Abra a câmera',
    'Foxbot, quantos banhos a Viih Tube toma por semana?',
    'Foxbot, por que foram comer o morcego?',
    'É necessário mexer o macarrão enquanto cozinha?',
    'Procure vídeos do Ronaldinho',
    'Como colocar uma camisinha?',
    'Ligue um timer por 15 minutos',
    'Quanto está o jogo do PSG?',
    'Qual é o próximo compromisso do dia?',
    'Foxbot, pare o cronômetro',
    'Por que o Thiago Leifert parece um sapatênis?',
    'O que fazer agora que acabou a água?',
    'Baixar tinta de impressora',
    'Tem alguma música ruim do Barões da Pisadinha?',
    'Amor ou o litrão?',
    'Quem é nóia é imune ao coronavírus?',
    'Tutorial de como enganar um agiota',
    'Como saber se ela gosta de mim?',
    'Como ficar rico com daytrade?',
    'O que é quarentena?',
    'Foxbot, meu colesterol está alto?',
    'Qual é o almoço de hoje?',
    'Quanto custa uma Itubaína',
    'Quantos anos vive uma lagosta?',
    'Por que liga da justiça é ruim?',
    'Mostre a nota de Corra no IMDB',
    'Como prender a respiração?',
    'Quando foi a última vez que tomei banho?',
    'Por que existem pessoas que riem como porcos?',
    'Onde encontro um dodo de estimação?',
    'O que são mulheres-sapiens?',
    'Postos de gasolina perto de mim',
    'Quero um milkshake de ovomaltine',
    'Por que curva de rio acumula lixo?',
    'Foxbot, pare o timer',
    'Mostre o meu próximo compromisso',
    'Toque minha playlist de funk',
    'O que é trap?',
    'Quanto tempo o cronômetro está marcando?',
    'Foxbot, me mostre uma receita de bolo',
    'Quanto está custando uma ação da Cyrela?',
    'Quem é Carlos Adão?',
    'Qual o nome do gorila que morreu em 2016?',
    'Como tratar queimaduras de sol?',
    'Qual a altura máxima de um anão?',
    'Diga a rota para o Maracanã',
    'Coronavac engravida?',
    'Como tá o câmbio de dólar pra real?',
    'Como saber se fui hackeado?',
    'Foxbot, marque dentista para amanhã às 10h',
    'Foxbot, ligue a impressora',
    'Foxbot, quantas calorias eu gastei hoje?',
    'Foxbot, toque uma música de ambiente',
    'Diminua o volume da música',
    'Foxbot, reserve uma mesa no Varanda para as 19h.',
    'Foxbot, acesse o site da Netshoes',
    'Como se fala eu amo você em russo?',
    'Foxbot, quando é o próximo show do Alok?',
    'Onde fica o Taj Mahal?',
    'Foxbot, quanto está uma passagem para Miami?',
    'Foxbot, quantas calorias tem uma torta de limão?',
    'Foxbot, mande uma mensagem para a minha irmã.'
)

$row = 201
foreach ($sentence in $sentences) {
    $ws.Range("A$row").Value = $intent
    $ws.Range("B$row").Value = $sentence
    $row = $row + 1
}

$ws.Range("B303").Select()
Write-Host "Added $($sentences.Count) sentences"
